# Update "想去人数" (column F) figures for the rows that changed in the
# latest scrape for both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 1065
    "F3"  = 759
    "F6"  = 1098
    "F8"  = 1793
    "F9"  = 6419
    "F10" = 487
    "F11" = 376
    "F16" = 6408
    "F17" = 277
    "F18" = 1295
    "F23" = 280
    "F24" = 111
    "F25" = 159
    "F29" = 394
    "F35" = 26
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
